$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 96
$srcRow = 95

# Write the new row's values first.
$ws.Cells.Item($row, 1).Value = 45447.2916666667
$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 1
$ws.Cells.Item($row, 5).Value = 1
$ws.Cells.Item($row, 6).Value = 1
$ws.Cells.Item($row, 8).Value = "YKY.MI"

# Column A needs the same date style (s="1") as the rest of the date column;
# copy formatting only from the row above so the existing style is reused
# instead of a brand-new one being created.
$ws.Cells.Item($srcRow, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column G stores "1" as text (shared string), not as a number. Copy the
# whole cell (value + format) from the identical cell above so it lands as
# the same shared string / plain style instead of being auto-typed as a
# number.
$ws.Cells.Item($srcRow, 7).Copy()
$ws.Cells.Item($row, 7).PasteSpecial(-4104)
$excel.CutCopyMode = $false
